$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.218619848223113
$ws.Range("C2").Value = 0.31505200208872
$ws.Range("E2").Value = 0.7874601494623903
$ws.Range("F2").Value = 2.040269078299843
$ws.Range("G2").Value = 0.002390410017571003
$ws.Range("I2").Value = 0.2888127484026981
$ws.Range("J2").Value = 0.02711371948510966
$ws.Range("O2").Value = 1.391930064394046
$ws.Range("B3").Value = 1.070181160103346
$ws.Range("C3").Value = 0.2751751077511813
$ws.Range("E3").Value = 0.7574822573187987
$ws.Range("F3").Value = 2.012894732451471
$ws.Range("G3").Value = 0.002393173187529599
$ws.Range("I3").Value = 0.2984102864338833
$ws.Range("J3").Value = 0.02766705527879587
$ws.Range("O3").Value = 1.419590983600742
$ws.Range("B4").Value = 0.9787304581741409
$ws.Range("C4").Value = 0.2505851961521728
$ws.Range("E4").Value = 0.7393691555516
$ws.Range("F4").Value = 1.99753969922665
$ws.Range("G4").Value = 0.002394956982225086
$ws.Range("I4").Value = 0.3046752336162761
$ws.Range("J4").Value = 0.02804103884288978
$ws.Range("O4").Value = 1.438319507967336
$ws.Range("B5").Value = 0.9413883067216489
$ws.Range("C5").Value = 0.2405387474883867
$ws.Range("E5").Value = 0.7320617851755316
$ws.Range("F5").Value = 1.991646645391739
$ws.Range("G5").Value = 0.002395705886897612
$ws.Range("I5").Value = 0.3073214411782148
$ws.Range("J5").Value = 0.02820200709854959
$ws.Range("O5").Value = 1.446388289314598
$ws.Range("B6").Value = 0.935183197058393
$ws.Range("C6").Value = 0.2388689990951036
$ws.Range("E6").Value = 0.7308528655964039
$ws.Range("F6").Value = 1.990690075790212
$ws.Range("G6").Value = 0.002395831571712855
$ws.Range("I6").Value = 0.3077664599698489
$ws.Range("J6").Value = 0.02822925196640114
$ws.Range("O6").Value = 1.447754426385856
$ws.Range("B7").Value = 0.9782271503476068
$ws.Range("C7").Value = 0.250449810098786
$ws.Range("E7").Value = 0.7392703066150688
$ws.Range("F7").Value = 1.997458750085841
$ws.Range("G7").Value = 0.002394966993357343
$ws.Range("I7").Value = 0.3047105444453977
$ws.Range("J7").Value = 0.02804317508456045
$ws.Range("O7").Value = 1.43842656054899
$ws.Range("B8").Value = 1.167503648461377
$ws.Range("C8").Value = 0.3013246741357989
$ws.Range("E8").Value = 0.7770629013912753
$ws.Range("F8").Value = 2.030528169412079
$ws.Range("G8").Value = 0.002391344701266055
$ws.Range("I8").Value = 0.2920445535007308
$ws.Range("J8").Value = 0.02729737978181745
$ws.Range("O8").Value = 1.401104291563158
$ws.Range("B9").Value = 1.536138342081699
$ws.Range("C9").Value = 0.4002317111374509
$ws.Range("E9").Value = 0.8535037947549995
$ws.Range("F9").Value = 2.106961771863539
$ws.Range("G9").Value = 0.002384930196905503
$ws.Range("I9").Value = 0.2701739534168772
$ws.Range("J9").Value = 0.0261083492357912
$ws.Range("O9").Value = 1.341844858408678
$ws.Range("B10").Value = 1.805340119875495
$ws.Range("C10").Value = 0.4723523277236836
$ws.Range("E10").Value = 0.9110940568714199
$ws.Range("F10").Value = 2.170267070673333
$ws.Range("G10").Value = 0.002380633039124305
$ws.Range("I10").Value = 0.2559353741355364
$ws.Range("J10").Value = 0.02540401904616729
$ws.Range("O10").Value = 1.30691763611064
$ws.Range("B11").Value = 1.927435620965412
$ws.Range("C11").Value = 0.5050388986937833
$ws.Range("E11").Value = 0.9376061060912662
$ws.Range("F11").Value = 2.200638219686113
$ws.Range("G11").Value = 0.002378767465089874
$ws.Range("I11").Value = 0.249859626009294
$ws.Range("J11").Value = 0.02512089743249923
$ws.Range("O11").Value = 1.29292358111941
$ws.Range("B12").Value = 1.973615485958476
$ws.Range("C12").Value = 0.5173984784631784
$ws.Range("E12").Value = 0.947690734644155
$ws.Range("F12").Value = 2.212366617459878
$ws.Range("G12").Value = 0.002378073781467414
$ws.Range("I12").Value = 0.2476170097586472
$ws.Range("J12").Value = 0.02501909284373305
$ws.Range("O12").Value = 1.287898929404733
$ws.Range("B13").Value = 1.96367231299439
$ws.Range("C13").Value = 0.5147374377523306
$ws.Range("E13").Value = 0.9455168248256314
$ws.Range("F13").Value = 2.209830558567006
$ws.Range("G13").Value = 0.002378222611753491
$ws.Range("I13").Value = 0.2480974050790721
$ws.Range("J13").Value = 0.0250407770480372
$ws.Range("O13").Value = 1.288968830987187
$ws.Range("B14").Value = 1.931235983530883
$ws.Range("C14").Value = 0.5060560955071196
$ws.Range("E14").Value = 0.9384348711550246
$ws.Range("F14").Value = 2.2015985552232
$ws.Range("G14").Value = 0.002378710139779695
$ws.Range("I14").Value = 0.2496739571078237
$ws.Range("J14").Value = 0.02511241326914337
$ws.Range("O14").Value = 1.292504684347037
$ws.Range("B15").Value = 1.911360553421162
$ws.Range("C15").Value = 0.5007361444378944
$ws.Range("E15").Value = 0.9341028406888086
$ws.Range("F15").Value = 2.196585880933128
$ws.Range("G15").Value = 0.002379010425876387
$ws.Range("I15").Value = 0.2506472241507636
$ws.Range("J15").Value = 0.02515699810633798
$ws.Range("O15").Value = 1.294706320991551
$ws.Range("B16").Value = 1.797353321214985
$ws.Range("C16").Value = 0.4702136836461364
$ws.Range("E16").Value = 0.9093677447669393
$ws.Range("F16").Value = 2.168314010961893
$ws.Range("G16").Value = 0.002380756748986625
$ws.Range("I16").Value = 0.2563405517926998
$ws.Range("J16").Value = 0.02542327587387483
$ws.Range("O16").Value = 1.307870488012
$ws.Range("B17").Value = 1.727318142711113
$ws.Range("C17").Value = 0.4514575690531615
$ws.Range("E17").Value = 0.8942739339029799
$ws.Range("F17").Value = 2.151373983596642
$ws.Range("G17").Value = 0.002381850870145938
$ws.Range("I17").Value = 0.2599363366065668
$ws.Range("J17").Value = 0.02559621036086845
$ws.Range("O17").Value = 1.316433048781107
$ws.Range("B18").Value = 1.687001480284039
$ws.Range("C18").Value = 0.4406581539925014
$ws.Range("E18").Value = 0.8856219264861949
$ws.Range("F18").Value = 2.141778533198405
$ws.Range("G18").Value = 0.002382488581713696
$ws.Range("I18").Value = 0.2620422990676747
$ws.Range("J18").Value = 0.02569918406096328
$ws.Range("O18").Value = 1.321536198486939
$ws.Range("B19").Value = 1.673345138942921
$ws.Range("C19").Value = 0.4369997206093785
$ws.Range("E19").Value = 0.8826975830232158
$ws.Range("F19").Value = 2.13855506003506
$ws.Range("G19").Value = 0.002382705944804716
$ws.Range("I19").Value = 0.2627618150000268
$ws.Range("J19").Value = 0.02573464999492003
$ws.Range("O19").Value = 1.323294575856877
$ws.Range("B20").Value = 1.734777072802331
$ws.Range("C20").Value = 0.4534553734438305
$ws.Range("E20").Value = 0.8958776375899618
$ws.Range("F20").Value = 2.153161951248222
$ws.Range("G20").Value = 0.002381733530071079
$ws.Range("I20").Value = 0.259549647970613
$ws.Range("J20").Value = 0.02557743795935075
$ws.Range("O20").Value = 1.315503091757037
$ws.Range("B21").Value = 1.940764833274329
$ws.Range("C21").Value = 0.5086065129415829
$ws.Range("E21").Value = 0.9405137892695734
$ws.Range("F21").Value = 2.204010309617416
$ws.Range("G21").Value = 0.002378566594902526
$ws.Range("I21").Value = 0.2492093044550225
$ws.Range("J21").Value = 0.02509122484912574
$ws.Range("O21").Value = 1.291458648877338
$ws.Range("B22").Value = 2.075067742792271
$ws.Range("C22").Value = 0.5445450857176866
$ws.Range("E22").Value = 0.9699488979020714
$ws.Range("F22").Value = 2.238569307771883
$ws.Range("G22").Value = 0.002376571217493167
$ws.Range("I22").Value = 0.2427903536795375
$ws.Range("J22").Value = 0.02480499663708358
$ws.Range("O22").Value = 1.277345758862253
$ws.Range("B23").Value = 2.003418019261744
$ws.Range("C23").Value = 0.5253738939756545
$ws.Range("E23").Value = 0.9542147964622956
$ws.Range("F23").Value = 2.220002708615766
$ws.Range("G23").Value = 0.002377629400842005
$ws.Range("I23").Value = 0.2461851089602067
$ws.Range("J23").Value = 0.02495486075554609
$ws.Range("O23").Value = 1.28473081755466
$ws.Range("B24").Value = 1.73140505387272
$ws.Range("C24").Value = 0.4525522167282361
$ws.Range("E24").Value = 0.895152523330438
$ws.Range("F24").Value = 2.152353163864092
$ws.Range("G24").Value = 0.002381786552443047
$ws.Range("I24").Value = 0.2597243493896899
$ws.Range("J24").Value = 0.02558591390364384
$ws.Range("O24").Value = 1.31592296348343
$ws.Range("B25").Value = 1.436693622159964
$ws.Range("C25").Value = 0.3735690541847703
$ws.Range("E25").Value = 0.8325741982617956
$ws.Range("F25").Value = 2.085035258337186
$ws.Range("G25").Value = 0.002386592198691313
$ws.Range("I25").Value = 0.2757707097392044
$ws.Range("J25").Value = 0.02640047472032592
$ws.Range("O25").Value = 1.356372568279085
